$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns (B: Coin, C: Link, D: Price, E: Volume) stay as text
# rather than being auto-coerced to numbers by Excel when they look numeric
# (e.g. "0.9969", "30.834.08", "36.00", "0.7300").

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.834.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.961.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.92%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9969"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6358"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +35.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9977"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.34%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3265"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +12.75%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +14.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06851"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8359"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +15.14%  "

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07979"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.02%  "

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "101.08"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.941.01"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.77%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.380"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.809.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007691"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.633"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.193.69"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9982"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9982"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.651"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.33%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.495"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.65%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.29"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.127"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1231"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +26.71%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.558"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.04%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.348"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.511"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.371"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05034"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.202"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.79%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7300"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01963"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.927"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.518"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.94%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.55"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.58%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4616"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.62%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.031"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.91%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8448"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9986"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.62"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.32%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.328"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4208"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "930.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.68%  "
